# Apply updated "想去人数" (interested-count) values scraped at commit 456a3b4.
# The workbook has two sheets that mirror the same event rows in column F:
#   "展览"   (Worksheets index/name based)
#   "全部类型"

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new value map for the "展览" sheet
$exhibitionUpdates = @{
    "F3"  = 13035
    "F6"  = 99
    "F10" = 13005
    "F11" = 293
    "F12" = 547
    "F13" = 8729
    "F14" = 7748
    "F15" = 208
    "F17" = 429
    "F18" = 132
    "F19" = 990
}

foreach ($cellRef in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range($cellRef).Value = $exhibitionUpdates[$cellRef]
}

# Row -> new value map for the "全部类型" sheet
$allTypesUpdates = @{
    "F4"  = 13035
    "F7"  = 99
    "F11" = 13005
    "F12" = 293
    "F13" = 547
    "F14" = 8729
    "F15" = 7748
    "F16" = 208
    "F18" = 429
    "F19" = 132
    "F20" = 990
}

foreach ($cellRef in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range($cellRef).Value = $allTypesUpdates[$cellRef]
}
